$d = $word.ActiveDocument

$replacements = @(
    @("2024-03-17 Sunday", "2024-03-18 Monday"),
    @("98×69=", "68×93="),
    @("80×83=", "22×97="),
    @("40×94=", "51×61="),
    @("96×95=", "19×44="),
    @("29×97=", "73×49="),
    @("15×14=", "62×45="),
    @("33×16=", "96×99="),
    @("40×60=", "60×35="),
    @("28×25=", "31×61="),
    @("11×88=", "34×77="),
    @("80×23=", "24×57="),
    @("21×27=", "46×35="),
    @("45×48=", "50×63="),
    @("39×15=", "41×92="),
    @("18×14=", "25×56="),
    @("80×39=", "73×29="),
    @("73×66=", "31×90="),
    @("73×63=", "90×27="),
    @("92×20=", "25×13="),
    @("36×98=", "18×88="),
    @("82×52=", "85×91="),
    @("68×23=", "25×87="),
    @("49×45=", "63×86="),
    @("59×32=", "69×57="),
    @("61×52=", "61×83=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
